# Updates the scraped "RPA Developer" job-listing data on the active sheet.
# The scraper now runs a full extra pass over the same search results
# (robot keeps working without stopping), so the original 10 data rows are
# refreshed with the latest posting info and then appended a second time
# as rows 12-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest scrape data: JobProfileName, CompanyName, DatePosted
$data = @(
    @("US Onsite Opportunity - Software Professionals", "DataQuest Corp", "7 days ago"),
    @("RPA Developer", "Volvo Group 4.1", "2 hours ago"),
    @("Associate 2 - RPA Developer - Bangalore", "KPMG 4.0", "30+ days ago"),
    @("UI Path & RPA Developer", "Menorah Personnel Management India Private", "20 days ago"),
    @("RPA Developer", "Infosys technology limited 3.9", "30+ days ago"),
    @("RPA Developer", "Krish Tech Inc", "25 days ago"),
    @("RPA Developer", "Capgemini 3.8", "30+ days ago"),
    @("Sr. RPA Developer", "Merck KGaA 4.0", "30+ days ago"),
    @("Developer - RPA", "Hudson's Bay Company 3.7", "30+ days ago"),
    @("Looking for RPA Developers", "Careator Technologies", "30+ days ago")
)

# Refresh rows 2-11 with the newest data
$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}

# The robot kept running and appended the same batch again as rows 12-21
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}
